# Auto-generated Excel COM-interop script applying numeric corrections
# to the Gilgamesh_Profits profit-calculation sheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR) per the updated Universalis price snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2381.95
$ws.Range("J17").Value = 2441.25
$ws.Range("L17").Value = 7323.75
$ws.Range("N17").Value = -7659.75
$ws.Range("H112").Value = 2113.1428
$ws.Range("J112").Value = 2113.1428
$ws.Range("L112").Value = 6339.428400000001
$ws.Range("N112").Value = -8555.428400000001
$ws.Range("H137").Value = 2542.6562
$ws.Range("I137").Value = 1903.8334
$ws.Range("J137").Value = 4459.125
$ws.Range("K137").Value = 5711.5002
$ws.Range("L137").Value = 13377.375
$ws.Range("M137").Value = -3161.5002
$ws.Range("N137").Value = -18477.375
$ws.Range("H138").Value = 585596.6
$ws.Range("J138").Value = 965030.1
$ws.Range("L138").Value = 2895090.3
$ws.Range("N138").Value = -2905370.3
$ws.Range("H141").Value = 2000
$ws.Range("I141").Value = 1000
$ws.Range("J141").Value = 2500
$ws.Range("K141").Value = 3000
$ws.Range("L141").Value = 7500
$ws.Range("M141").Value = 2180
$ws.Range("N141").Value = -17860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2643.89
$ws.Range("I32").Value = 2267.0532
$ws.Range("K32").Value = 2267.0532
$ws.Range("M32").Value = -1980.0532
$ws.Range("H45").Value = 26529.111
$ws.Range("I45").Value = 62267.855
$ws.Range("J45").Value = 3786.2727
$ws.Range("K45").Value = 62267.855
$ws.Range("L45").Value = 3786.2727
$ws.Range("M45").Value = -61890.855
$ws.Range("N45").Value = -4540.2727
$ws.Range("H74").Value = 161888.11
$ws.Range("I74").Value = 244696.39
$ws.Range("J74").Value = 3172.25
$ws.Range("K74").Value = 244696.39
$ws.Range("L74").Value = 3172.25
$ws.Range("M74").Value = -243822.39
$ws.Range("N74").Value = -4920.25
$ws.Range("H77").Value = 161888.11
$ws.Range("I77").Value = 244696.39
$ws.Range("J77").Value = 3172.25
$ws.Range("K77").Value = 1223481.95
$ws.Range("L77").Value = 15861.25
$ws.Range("M77").Value = -1219113.95
$ws.Range("N77").Value = -24597.25
$ws.Range("H122").Value = 5080.727
$ws.Range("I122").Value = 4125.4375
$ws.Range("K122").Value = 12376.3125
$ws.Range("M122").Value = -9926.3125
$ws.Range("H132").Value = 2054.0784
$ws.Range("I132").Value = 1611.317
$ws.Range("K132").Value = 4833.951
$ws.Range("M132").Value = -2303.951

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 40000796
$ws.Range("I94").Value = 52632136
$ws.Range("K94").Value = 52632136
$ws.Range("M94").Value = -52631685
$ws.Range("H132").Value = 83999.75
$ws.Range("J132").Value = 83999.75
$ws.Range("L132").Value = 83999.75
$ws.Range("N132").Value = -94119.75
$ws.Range("H134").Value = 2648.7144
$ws.Range("I134").Value = 1883.25
$ws.Range("K134").Value = 5649.75
$ws.Range("M134").Value = -3114.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3535.25
$ws.Range("I31").Value = 2811.8655
$ws.Range("J31").Value = 5416.05
$ws.Range("K31").Value = 2811.8655
$ws.Range("L31").Value = 5416.05
$ws.Range("M31").Value = -2516.8655
$ws.Range("N31").Value = -6006.05
$ws.Range("H34").Value = 3535.25
$ws.Range("I34").Value = 2811.8655
$ws.Range("J34").Value = 5416.05
$ws.Range("K34").Value = 2811.8655
$ws.Range("L34").Value = 5416.05
$ws.Range("M34").Value = -2609.8655
$ws.Range("N34").Value = -5820.05
$ws.Range("H69").Value = 4086
$ws.Range("I69").Value = 4086
$ws.Range("K69").Value = 4086
$ws.Range("M69").Value = -3337
$ws.Range("H72").Value = 4086
$ws.Range("I72").Value = 4086
$ws.Range("K72").Value = 12258
$ws.Range("M72").Value = -8514
$ws.Range("H99").Value = 17003.834
$ws.Range("I99").Value = 17003.834
$ws.Range("K99").Value = 17003.834
$ws.Range("M99").Value = -15505.834
$ws.Range("H126").Value = 17003.834
$ws.Range("I126").Value = 17003.834
$ws.Range("K126").Value = 51011.50199999999
$ws.Range("M126").Value = -48541.50199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 8595.200000000001
$ws.Range("I80").Value = 8493.333000000001
$ws.Range("K80").Value = 25479.999
$ws.Range("M80").Value = -24543.999
$ws.Range("H83").Value = 8595.200000000001
$ws.Range("I83").Value = 8493.333000000001
$ws.Range("K83").Value = 76439.997
$ws.Range("M83").Value = -71759.997
$ws.Range("H122").Value = 2354.75
$ws.Range("J122").Value = 2471.6667
$ws.Range("L122").Value = 22245.0003
$ws.Range("N122").Value = -27145.0003
$ws.Range("H132").Value = 3421.64
$ws.Range("J132").Value = 3632.2
$ws.Range("L132").Value = 32689.8
$ws.Range("N132").Value = -37749.8
$ws.Range("H138").Value = 7958.25
$ws.Range("I138").Value = 9500
$ws.Range("K138").Value = 28500
$ws.Range("M138").Value = -23360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 83336480
$ws.Range("I80").Value = 125002540
$ws.Range("J80").Value = 4351
$ws.Range("K80").Value = 125002540
$ws.Range("L80").Value = 4351
$ws.Range("M80").Value = -125001542
$ws.Range("N80").Value = -6347
$ws.Range("H83").Value = 83336480
$ws.Range("I83").Value = 125002540
$ws.Range("J83").Value = 4351
$ws.Range("K83").Value = 625012700
$ws.Range("L83").Value = 21755
$ws.Range("M83").Value = -625007708
$ws.Range("N83").Value = -31739
$ws.Range("H122").Value = 2119.5454
$ws.Range("I122").Value = 2111.7
$ws.Range("J122").Value = 2198
$ws.Range("K122").Value = 6335.099999999999
$ws.Range("L122").Value = 6594
$ws.Range("M122").Value = -3885.099999999999
$ws.Range("N122").Value = -11494
$ws.Range("H132").Value = 3356.3914
$ws.Range("I132").Value = 2816
$ws.Range("K132").Value = 8448
$ws.Range("M132").Value = -5918
$ws.Range("H134").Value = 62527.715
$ws.Range("J134").Value = 62527.715
$ws.Range("L134").Value = 187583.145
$ws.Range("N134").Value = -192653.145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I2").Value = 9999
$ws.Range("K2").Value = 9999
$ws.Range("M2").Value = -9887
$ws.Range("H22").Value = 2026.5
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 2026.5
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H40").Value = 21201.225
$ws.Range("I40").Value = 34613.816
$ws.Range("K40").Value = 34613.816
$ws.Range("M40").Value = -34477.816
$ws.Range("H53").Value = 49950
$ws.Range("J53").Value = 49950
$ws.Range("L53").Value = 49950
$ws.Range("N53").Value = -50986
$ws.Range("H55").Value = 639.6667
$ws.Range("I55").Value = 607.6
$ws.Range("J55").Value = 800
$ws.Range("K55").Value = 607.6
$ws.Range("L55").Value = 800
$ws.Range("M55").Value = -434.6
$ws.Range("N55").Value = -1146
$ws.Range("H93").Value = 1529.1111
$ws.Range("I93").Value = 1529.1111
$ws.Range("K93").Value = 1529.1111
$ws.Range("M93").Value = -281.1111000000001
$ws.Range("H100").Value = 3724.2334
$ws.Range("I100").Value = 3210.48
$ws.Range("J100").Value = 6293
$ws.Range("K100").Value = 3210.48
$ws.Range("L100").Value = 6293
$ws.Range("M100").Value = -2669.48
$ws.Range("N100").Value = -7375
$ws.Range("H132").Value = 4049.625
$ws.Range("I132").Value = 3399.5
$ws.Range("K132").Value = 10198.5
$ws.Range("M132").Value = -7668.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 3999.8572
$ws.Range("J14").Value = 4999.5
$ws.Range("L14").Value = 4999.5
$ws.Range("N14").Value = -5335.5
$ws.Range("H46").Value = 64915.57
$ws.Range("J46").Value = 64915.57
$ws.Range("L46").Value = 64915.57
$ws.Range("N46").Value = -65377.57
$ws.Range("H53").Value = 59188
$ws.Range("I53").Value = 59188
$ws.Range("K53").Value = 59188
$ws.Range("M53").Value = -58581
$ws.Range("H58").Value = 14874.25
$ws.Range("I58").Value = 4999.5
$ws.Range("K58").Value = 4999.5
$ws.Range("M58").Value = -4691.5
$ws.Range("H96").Value = 4551.364
$ws.Range("I96").Value = 4353
$ws.Range("J96").Value = 4898.5
$ws.Range("K96").Value = 4353
$ws.Range("L96").Value = 4898.5
$ws.Range("M96").Value = -2980
$ws.Range("N96").Value = -7644.5
$ws.Range("H122").Value = 8930919
$ws.Range("I122").Value = 1988.5416
$ws.Range("J122").Value = 62504504
$ws.Range("K122").Value = 5965.6248
$ws.Range("L122").Value = 187513512
$ws.Range("M122").Value = -3515.6248
$ws.Range("N122").Value = -187518412
$ws.Range("H126").Value = 1839.6
$ws.Range("I126").Value = 1747.7778
$ws.Range("K126").Value = 5243.3334
$ws.Range("M126").Value = -2773.3334
$ws.Range("H132").Value = 4177.696
$ws.Range("I132").Value = 4378.121
$ws.Range("K132").Value = 13134.363
$ws.Range("M132").Value = -10604.363
$ws.Range("H134").Value = 64915.57
$ws.Range("J134").Value = 64915.57
$ws.Range("L134").Value = 194746.71
$ws.Range("N134").Value = -199816.71
